$d = $word.ActiveDocument

# 1) "Poate vedea și căuta produsele existente deja pe site"
#    -> "Poate vedea produsele existente deja pe site"
$d.Content.Find.Execute(
    "Poate vedea și căuta produsele existente deja pe site",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Poate vedea produsele existente deja pe site", 2)

# 2) Add a new bullet after "Poate adăuga în coș produsele dorite pentru cumpărare"
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Poate face o comandă și poate să o plătească online"
